# "Add ARM C and D" - update the Bento Filter Arm-B startup sheet:
#   - rename the generated TSV/Web companion workbook filenames from
#     TC01_... to TC02_... (shared-string text content), and
#   - move the saved cell selection from E3 to D3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("startup")

# Update the two filename cells (D2/E2) that reference the companion
# TC01_Bento_Filter_Arm-B_*.xlsx outputs -> TC02_Bento_Filter_Arm-B_*.xlsx
$ws.Range("D2").Value = "TC02_Bento_Filter_Arm-B_TSVData.xlsx"
$ws.Range("E2").Value = "TC02_Bento_Filter_Arm-B_WebData.xlsx"

# Move the active selection on the sheet from E3 to D3
$ws.Range("D3").Select() | Out-Null
